# Scheduled-runner style market-price refresh across all Leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Updates the per-leve columns
# H:N (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# with freshly retrieved marketboard figures. Plain data values only -
# no formulas are used anywhere in this workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 15000
$ws.Range("J16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("N16").Value = -15460

# Row 28
$ws.Range("H28").Value = 745.6
$ws.Range("I28").Value = 435.7143
$ws.Range("K28").Value = 435.7143
$ws.Range("M28").Value = 49.28570000000002

# Row 51
$ws.Range("H51").Value = 7666.6665
$ws.Range("J51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("N51").Value = -8968

# Row 95
$ws.Range("H95").Value = 6000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 6000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 6000
$ws.Range("M95").Value = $null
$ws.Range("N95").Value = -11492

# Row 103
$ws.Range("H103").Value = 417051.34
$ws.Range("J103").Value = 433.33334
$ws.Range("L103").Value = 1300.00002
$ws.Range("N103").Value = -2472.00002

# Row 132
$ws.Range("H132").Value = 3286.6667
$ws.Range("I132").Value = 3799.5
$ws.Range("J132").Value = 1235.3334
$ws.Range("K132").Value = 11398.5
$ws.Range("L132").Value = 3706.0002
$ws.Range("M132").Value = -8868.5
$ws.Range("N132").Value = -8766.0002

# Row 137
$ws.Range("H137").Value = 1712.4166
$ws.Range("I137").Value = 1699.8572
$ws.Range("J137").Value = 1730
$ws.Range("K137").Value = 5099.571599999999
$ws.Range("L137").Value = 5190
$ws.Range("M137").Value = -2549.571599999999
$ws.Range("N137").Value = -10290

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 109.75
$ws.Range("I5").Value = 39
$ws.Range("J5").Value = 133.33333
$ws.Range("K5").Value = 39
$ws.Range("L5").Value = 133.33333
$ws.Range("M5").Value = 73
$ws.Range("N5").Value = -357.33333

# Row 45
$ws.Range("H45").Value = 2345.4546
$ws.Range("I45").Value = 1665.8667
$ws.Range("J45").Value = 3801.7144
$ws.Range("K45").Value = 1665.8667
$ws.Range("L45").Value = 3801.7144
$ws.Range("M45").Value = -1288.8667
$ws.Range("N45").Value = -4555.7144

# Row 57
$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 3000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2516

# Row 61
$ws.Range("H61").Value = 1615.3077
$ws.Range("I61").Value = 1537.4166
$ws.Range("K61").Value = 1537.4166
$ws.Range("M61").Value = -1325.4166

# Row 92
$ws.Range("H92").Value = 20650
$ws.Range("J92").Value = 20650
$ws.Range("L92").Value = 20650
$ws.Range("N92").Value = -25642

# Row 102
$ws.Range("H102").Value = 1319.625
$ws.Range("I102").Value = 1242.8334
$ws.Range("J102").Value = 1550
$ws.Range("K102").Value = 1242.8334
$ws.Range("L102").Value = 1550
$ws.Range("M102").Value = 379.1666
$ws.Range("N102").Value = -4794

# Row 136
$ws.Range("H136").Value = 1615.3077
$ws.Range("I136").Value = 1537.4166
$ws.Range("K136").Value = 4612.2498
$ws.Range("M136").Value = -2062.2498

# Row 139
$ws.Range("H139").Value = 37547.637
$ws.Range("J139").Value = 37547.637
$ws.Range("L139").Value = 37547.637
$ws.Range("N139").Value = -47827.637

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 109.75
$ws.Range("I4").Value = 39
$ws.Range("J4").Value = 133.33333
$ws.Range("K4").Value = 39
$ws.Range("L4").Value = 133.33333
$ws.Range("M4").Value = 76
$ws.Range("N4").Value = -363.33333

# Row 22
$ws.Range("H22").Value = 519.9474
$ws.Range("I22").Value = 454.9375
$ws.Range("K22").Value = 454.9375
$ws.Range("M22").Value = -281.9375

# Row 64
$ws.Range("H64").Value = 832.3333
$ws.Range("I64").Value = 1117.5
$ws.Range("J64").Value = 475.875
$ws.Range("K64").Value = 1117.5
$ws.Range("L64").Value = 475.875
$ws.Range("M64").Value = -892.5
$ws.Range("N64").Value = -925.875

# Row 67
$ws.Range("H67").Value = 832.3333
$ws.Range("I67").Value = 1117.5
$ws.Range("J67").Value = 475.875
$ws.Range("K67").Value = 1117.5
$ws.Range("L67").Value = 475.875
$ws.Range("M67").Value = -337.5
$ws.Range("N67").Value = -2035.875

# Row 81
$ws.Range("H81").Value = 12908.546
$ws.Range("J81").Value = 12908.546
$ws.Range("L81").Value = 12908.546
$ws.Range("N81").Value = -15030.546

# Row 84
$ws.Range("H84").Value = 12908.546
$ws.Range("J84").Value = 12908.546
$ws.Range("L84").Value = 38725.638
$ws.Range("N84").Value = -49333.638

# Row 99
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

# Row 107
$ws.Range("H107").Value = 646.7857
$ws.Range("I107").Value = 679.75
$ws.Range("J107").Value = 449
$ws.Range("K107").Value = 679.75
$ws.Range("L107").Value = 449
$ws.Range("M107").Value = 1240.25
$ws.Range("N107").Value = -4289

# Row 135
$ws.Range("H135").Value = 42267.6
$ws.Range("J135").Value = 42267.6
$ws.Range("L135").Value = 42267.6
$ws.Range("N135").Value = -52407.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 32327.666
$ws.Range("J86").Value = 49355.332
$ws.Range("L86").Value = 49355.332
$ws.Range("N86").Value = -51601.332

# Row 89
$ws.Range("H89").Value = 32327.666
$ws.Range("J89").Value = 49355.332
$ws.Range("L89").Value = 246776.66
$ws.Range("N89").Value = -258008.66

# Row 105
$ws.Range("H105").Value = 804.84
$ws.Range("I105").Value = 738.63635
$ws.Range("J105").Value = 1290.3334
$ws.Range("K105").Value = 738.63635
$ws.Range("L105").Value = 1290.3334
$ws.Range("M105").Value = 1008.36365
$ws.Range("N105").Value = -4784.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 36.666668
$ws.Range("I2").Value = 19.333334
$ws.Range("J2").Value = 123.333336
$ws.Range("K2").Value = 116.000004
$ws.Range("L2").Value = 740.000016
$ws.Range("M2").Value = -3.000004000000004
$ws.Range("N2").Value = -966.000016

# Row 12
$ws.Range("H12").Value = 69.42856999999999
$ws.Range("I12").Value = 4.8333335
$ws.Range("J12").Value = 117.875
$ws.Range("K12").Value = 14.5000005
$ws.Range("L12").Value = 353.625
$ws.Range("M12").Value = 158.4999995
$ws.Range("N12").Value = -699.625

# Row 60
$ws.Range("H60").Value = 140
$ws.Range("I60").Value = 140
$ws.Range("K60").Value = 420
$ws.Range("M60").Value = -169

# Row 122
$ws.Range("H122").Value = 733.5294
$ws.Range("I122").Value = 399.5
$ws.Range("K122").Value = 3595.5
$ws.Range("M122").Value = -1145.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = $null

# Row 97
$ws.Range("H97").Value = 1220.95
$ws.Range("I97").Value = 1112.8823
$ws.Range("J97").Value = 1833.3334
$ws.Range("K97").Value = 1112.8823
$ws.Range("L97").Value = 1833.3334
$ws.Range("M97").Value = -616.8823
$ws.Range("N97").Value = -2825.3334

# Row 107
$ws.Range("H107").Value = 2849196.8
$ws.Range("I107").Value = 182.88235
$ws.Range("K107").Value = 182.88235
$ws.Range("M107").Value = 1737.11765

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2498
$ws.Range("I93").Value = 2495
$ws.Range("K93").Value = 2495
$ws.Range("M93").Value = -1247

# Row 127
$ws.Range("H127").Value = 39918.688
$ws.Range("J127").Value = 39918.688
$ws.Range("L127").Value = 39918.688
$ws.Range("N127").Value = -49838.688

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

# Row 113
$ws.Range("H113").Value = 1335.8125
$ws.Range("J113").Value = 382.44446
$ws.Range("L113").Value = 1147.33338
$ws.Range("N113").Value = -5487.33338

# Row 126
$ws.Range("H126").Value = 1230.5
$ws.Range("I126").Value = 1251.0476
$ws.Range("J126").Value = 799
$ws.Range("K126").Value = 3753.142800000001
$ws.Range("L126").Value = 2397
$ws.Range("M126").Value = -1283.142800000001
$ws.Range("N126").Value = -7337
